# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeout) values for rows 2..55, replacing the old "Strike#" values
# that were previously stored in column G.
$kValues = @(1, 2, 3, 0, 0, 1, 0, 2, 1, 2, 1, 1, 3, 0, 1, 1, 1, 0, 3, 3, 1, 2, 2, 1, 2, 1, 1, 3, 0, 2, 1, 1, 1, 1, 1, 1, 0, 0, 1, 0, 3, 2, 1, 0, 1, 0, 1, 1, 1, 2, 1, 0, 1, 2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
